# Apply 3 dpi phage count data for plate experiment 2
# - Updates PFU counts in column D for rows 98-145, 170-175, 206-211
#   (some rows previously had a placeholder 1000x-too-large value, others were
#    blank cells awaiting data, and some blank trailing-zero cells simply
#    needed the scientific-notation number format applied).
# - Updates the saved sheet view (active cell / selection) to D146.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dColumnValues = @{
    98 = "120000"
    99 = "320000"
    100 = "2000000"
    101 = "180000"
    102 = "160000"
    103 = "460000"
    110 = "5800000000"
    111 = "0"
    112 = "3.6E-9"
    113 = "0"
    114 = "46000000000"
    115 = "0"
    116 = "700000000"
    117 = "1000000000000"
    118 = "26000000000"
    119 = "3400000000"
    120 = "18000000000"
    121 = "0"
    122 = "30000000"
    123 = "2600000000"
    124 = "520000000000"
    125 = "1200000000000"
    126 = "400"
    127 = "26000000000"
    128 = "260000"
    129 = "3200000"
    130 = "1000"
    131 = "3200000000"
    132 = "52000000"
    133 = "38000000"
    134 = "0"
    135 = "0"
    136 = "0"
    137 = "1200"
    138 = "360"
    139 = "0"
    140 = "0"
    141 = "0"
    142 = "600"
    143 = "0"
    144 = "200"
    145 = "0"
    170 = "0"
    171 = "0"
    172 = "0"
    173 = "0"
    174 = "0"
    175 = "0"
    206 = "0"
    207 = "0"
    208 = "0"
    209 = "0"
    210 = "0"
    211 = "0"
}

foreach ($row in $dColumnValues.Keys) {
    $cell = $ws.Cells.Item([int]$row, 4)
    $cell.Value = $dColumnValues[$row]
    $cell.NumberFormat = "0.00E+00"
}

# Update the view: scroll position / active selection moved to D146
[void]$ws.Range("D146").Select()
